# Update Mayank Agarwal innings stats (runs/balls/fours/sixes) for rows 2-12.
# Values are stored as text (matching the workbook's existing "number stored
# as text" convention), so the number format is forced to text before the
# values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the individual cells whose value actually changes are listed here.
# (Row 10, and the F column on rows 3/4/6/12, keep their original values and
# are intentionally left untouched.)
$updates = @(
    @{ Cell = "C2";  Value = "45" }
    @{ Cell = "D2";  Value = "25" }
    @{ Cell = "E2";  Value = "4" }
    @{ Cell = "F2";  Value = "3" }

    @{ Cell = "C3";  Value = "11" }
    @{ Cell = "D3";  Value = "10" }
    @{ Cell = "E3";  Value = "1" }

    @{ Cell = "C4";  Value = "5" }
    @{ Cell = "D4";  Value = "9" }
    @{ Cell = "E4";  Value = "0" }

    @{ Cell = "C5";  Value = "26" }
    @{ Cell = "D5";  Value = "15" }
    @{ Cell = "E5";  Value = "5" }
    @{ Cell = "F5";  Value = "0" }

    @{ Cell = "C6";  Value = "9" }
    @{ Cell = "D6";  Value = "6" }
    @{ Cell = "E6";  Value = "1" }

    @{ Cell = "C7";  Value = "56" }
    @{ Cell = "D7";  Value = "39" }
    @{ Cell = "E7";  Value = "6" }
    @{ Cell = "F7";  Value = "1" }

    @{ Cell = "C8";  Value = "26" }
    @{ Cell = "D8";  Value = "19" }
    @{ Cell = "E8";  Value = "3" }
    @{ Cell = "F8";  Value = "0" }

    @{ Cell = "C9";  Value = "89" }
    @{ Cell = "D9";  Value = "60" }
    @{ Cell = "E9";  Value = "7" }
    @{ Cell = "F9";  Value = "4" }

    @{ Cell = "C11"; Value = "106" }
    @{ Cell = "D11"; Value = "50" }
    @{ Cell = "E11"; Value = "10" }
    @{ Cell = "F11"; Value = "7" }

    @{ Cell = "C12"; Value = "26" }
    @{ Cell = "D12"; Value = "20" }
    @{ Cell = "E12"; Value = "4" }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
}
